$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new export timestamp
$ws.Name = "IClientBalance-20240806-075647-"

# Column G (rows 2-274) holds a date serial that needs to move forward by one day
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
